$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the swapped variance terms feeding the a_0 / a_1 confidence
#    intervals: the adjugate-matrix block (C21:D22) must be used instead of
#    the original P matrix (C10:D11).
# ---------------------------------------------------------------------------
$ws.Range("C25:D26").FormulaArray = "=C21:D22*D17"

# ---------------------------------------------------------------------------
# 2. Confidence bounds for the regression function now use the single
#    critical value in $R$32 for every column (both bounds previously,
#    mistakenly, used S32 for columns AB:AF, which is empty).
# ---------------------------------------------------------------------------
$ws.Range("AA52").Formula = '=AA18-$R$32*SQRT(AA38)'
$ws.Range("AB52:AF52").Formula = '=AB18-$R$32*SQRT(AB38)'

$ws.Range("AA53").Formula = '=AA18+$R$32*SQRT(AA38)'
$ws.Range("AB53:AF53").Formula = '=AB18+$R$32*SQRT(AB38)'

# ---------------------------------------------------------------------------
# 3. Rework the printed a_0 / a_1 confidence-interval summary block
#    (rows 43-48): row 43 keeps the symbolic a_0 bounds, row 44 carries the
#    numeric a_0 bounds, row 47 now carries the symbolic a_1 bounds (moved
#    down from row 46) and row 48 (new) carries the numeric a_1 bounds
#    (moved down from row 47).
# ---------------------------------------------------------------------------

# Row 44: numeric a_0 interval: 16.038 < a_0 < 26.05
$ws.Range("AE44").Value = 16.038
$ws.Range("AF44").Value = "< a_0 <"
$ws.Range("AG44").Value = 26.05

# Row 46: clear the old a_1 symbolic-bounds row (moved to row 47)
$ws.Range("AC46").ClearContents()
$ws.Range("AD46").Clear()
$ws.Range("AF46").Clear()
$ws.Range("AG46").Clear()

# Row 47: clear the old numeric a_1 bounds (moved to row 48) and place the
# symbolic a_1 bounds here instead.
$ws.Range("AE47").Clear()
$ws.Range("AD47").Value = "'21,044-2,78*sqrt{0,007}"
$ws.Range("AF47").Value = "< a_1 <"
$ws.Range("AG47").Value = "'21,044+2,78*sqrt{0.007}"

# Row 48 (new): numeric a_1 interval: -2.104 < a_1 < -1.638
$ws.Range("AE48").Value = -2.104
$ws.Range("AF48").Value = "< a_1 <"
$ws.Range("AG48").Value = -1.638

# ---------------------------------------------------------------------------
# 4. Drop the stray formatted-only block to the right of the charts
#    (AH50:AO55) - no longer needed.
# ---------------------------------------------------------------------------
$ws.Range("AH50:AO55").Clear()

# ---------------------------------------------------------------------------
# 5. Restore the view to where the author left it.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("Y51").Select()
